$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" '42.197.65'
Set-TextValue "E2" '  -1.06%  '

Set-TextValue "D3" '2.265.82'
Set-TextValue "E3" '  -1.41%  '

Set-TextValue "E4" '  -0.10%  '

Set-TextValue "D5" '306.33'
Set-TextValue "E5" '  -0.50%  '

Set-TextValue "D6" '96.73'
Set-TextValue "E6" '  -0.03%  '

Set-TextValue "D7" '0.528'
Set-TextValue "E7" '  -1.23%  '

Set-TextValue "E8" '  -0.05%  '

Set-TextValue "E9" '  -1.31%  '

Set-TextValue "D10" '35.02'
Set-TextValue "E10" '  -2.06%  '

Set-TextValue "D11" '0.0790'
Set-TextValue "E11" '  -2.46%  '

Set-TextValue "E12" '  -0.19%  '

Set-TextValue "D13" '6.93'
Set-TextValue "E13" '  +2.39%  '

Set-TextValue "D14" '2.617.79'
Set-TextValue "E14" '  -1.44%  '

Set-TextValue "D15" '14.69'
Set-TextValue "E15" '  +0.86%  '

Set-TextValue "D16" '2.269.78'
Set-TextValue "E16" '  -1.90%  '

Set-TextValue "E17" '  -1.44%  '

Set-TextValue "D18" '42.080.70'
Set-TextValue "E18" '  -1.16%  '

Set-TextValue "D19" '12.32'
Set-TextValue "E19" '  -3.85%  '

Set-TextValue "D20" '0.0₃0905'
Set-TextValue "E20" '  -1.96%  '

Set-TextValue "D21" '6.01'
Set-TextValue "E21" '  -0.60%  '

Set-TextValue "D22" '67.79'
Set-TextValue "E22" '  -0.84%  '

Set-TextValue "D23" '237.20'
Set-TextValue "E23" '  -2.94%  '

Set-TextValue "D24" '2.57'
Set-TextValue "E24" '  -1.92%  '

Set-TextValue "D25" '1.96'
Set-TextValue "E25" '  -0.73%  '

Set-TextValue "E26" '  +0.01%  '

Set-TextValue "D27" '23.52'
Set-TextValue "E27" '  -3.23%  '

Set-TextValue "D28" '37.77'
Set-TextValue "E28" '  +2.40%  '

Set-TextValue "D29" '9.58'
Set-TextValue "E29" '  -1.41%  '

Set-TextValue "E30" '  +0.07%  '

Set-TextValue "D31" '163.15'
Set-TextValue "E31" '  +0.88%  '

Set-TextValue "D32" '5.24'
Set-TextValue "E32" '  -2.70%  '

Set-TextValue "E33" '  -0.06%  '

Set-TextValue "E34" '  +1.67%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D35" '0.0738'
Set-TextValue "E35" '  -2.83%  '

$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue "D36" '17.60'
Set-TextValue "E36" '  +1.31%  '

Set-TextValue "E37" '  -1.00%  '

Set-TextValue "D38" '0.104'
Set-TextValue "E38" '  -4.50%  '

Set-TextValue "E39" '  -1.38%  '

Set-TextValue "E40" '  -1.80%  '

Set-TextValue "D41" '4.05'
Set-TextValue "E41" '  -4.01%  '

Set-TextValue "E42" '  +2.73%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D43" '1.949.86'
Set-TextValue "E43" '  -3.69%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D44" '18.96'
Set-TextValue "E44" '  -4.36%  '

Set-TextValue "E45" '  -1.63%  '

Set-TextValue "E46" '  -3.64%  '

Set-TextValue "D47" '9.85'
Set-TextValue "E47" '  -4.31%  '

Set-TextValue "D48" '53.94'
Set-TextValue "E48" '  -0.04%  '

Set-TextValue "D49" '92.23'
Set-TextValue "E49" '  -1.07%  '

Set-TextValue "D50" '71.66'
Set-TextValue "E50" '  -2.55%  '

Set-TextValue "E51" '  -2.75%  '
